$wb = $excel.ActiveWorkbook

# --- "2o Parcial" sheet updates ---
$wsSegundo = $wb.Worksheets.Item("2o Parcial")
$wsSegundo.Range("E2").Value = 10
$wsSegundo.Range("F2").Value = 0
$wsSegundo.Range("G2").Value = 100
$wsSegundo.Range("H2").Value = 0
$wsSegundo.Range("I2").Value = 8.5
$wsSegundo.Range("J2").Value = 0
$wsSegundo.Range("K2").Value = 0
$wsSegundo.Range("E3").Value = 10
$wsSegundo.Range("F3").Value = 0
$wsSegundo.Range("G3").Value = 100
$wsSegundo.Range("H3").Value = 0
$wsSegundo.Range("I3").Value = 8.5
$wsSegundo.Range("J3").Value = 0
$wsSegundo.Range("K3").Value = 0
$wsSegundo.Range("E6").Value = 29
$wsSegundo.Range("F6").Value = 8
$wsSegundo.Range("G6").Value = 78.40000000000001
$wsSegundo.Range("H6").Value = 21.6
$wsSegundo.Range("I6").Value = 7.7
$wsSegundo.Range("J6").Value = 0
$wsSegundo.Range("K6").Value = 0
$wsSegundo.Range("E7").Value = 39
$wsSegundo.Range("F7").Value = 0
$wsSegundo.Range("G7").Value = 100
$wsSegundo.Range("H7").Value = 0
$wsSegundo.Range("I7").Value = 8.4
$wsSegundo.Range("J7").Value = 0
$wsSegundo.Range("K7").Value = 0
$wsSegundo.Range("E8").Value = 27
$wsSegundo.Range("F8").Value = 5
$wsSegundo.Range("G8").Value = 84.40000000000001
$wsSegundo.Range("H8").Value = 15.6
$wsSegundo.Range("I8").Value = 7.8
$wsSegundo.Range("J8").Value = 0
$wsSegundo.Range("K8").Value = 0
$wsSegundo.Range("E9").Value = 35
$wsSegundo.Range("F9").Value = 4
$wsSegundo.Range("G9").Value = 89.7
$wsSegundo.Range("H9").Value = 10.3
$wsSegundo.Range("I9").Value = 7.1
$wsSegundo.Range("J9").Value = 0
$wsSegundo.Range("K9").Value = 0
$wsSegundo.Range("E10").Value = 32
$wsSegundo.Range("F10").Value = 6
$wsSegundo.Range("G10").Value = 84.2
$wsSegundo.Range("H10").Value = 15.8
$wsSegundo.Range("I10").Value = 7.7
$wsSegundo.Range("J10").Value = 0
$wsSegundo.Range("K10").Value = 0
$wsSegundo.Range("E11").Value = 21
$wsSegundo.Range("F11").Value = 3
$wsSegundo.Range("G11").Value = 87.5
$wsSegundo.Range("H11").Value = 12.5
$wsSegundo.Range("I11").Value = 7.5
$wsSegundo.Range("J11").Value = 0
$wsSegundo.Range("K11").Value = 0
$wsSegundo.Range("E12").Value = 34
$wsSegundo.Range("F12").Value = 0
$wsSegundo.Range("G12").Value = 100
$wsSegundo.Range("H12").Value = 0
$wsSegundo.Range("I12").Value = 7.9
$wsSegundo.Range("J12").Value = 0
$wsSegundo.Range("K12").Value = 0
$wsSegundo.Range("E13").Value = 31
$wsSegundo.Range("F13").Value = 0
$wsSegundo.Range("G13").Value = 100
$wsSegundo.Range("H13").Value = 0
$wsSegundo.Range("I13").Value = 8
$wsSegundo.Range("J13").Value = 0
$wsSegundo.Range("K13").Value = 0
$wsSegundo.Range("E14").Value = 39
$wsSegundo.Range("F14").Value = 1
$wsSegundo.Range("G14").Value = 97.5
$wsSegundo.Range("H14").Value = 2.5
$wsSegundo.Range("I14").Value = 8.9
$wsSegundo.Range("J14").Value = 0
$wsSegundo.Range("K14").Value = 0
$wsSegundo.Range("E15").Value = 24
$wsSegundo.Range("F15").Value = 0
$wsSegundo.Range("G15").Value = 100
$wsSegundo.Range("H15").Value = 0
$wsSegundo.Range("I15").Value = 7.6
$wsSegundo.Range("J15").Value = 0
$wsSegundo.Range("K15").Value = 0
$wsSegundo.Range("E16").Value = 35
$wsSegundo.Range("F16").Value = 3
$wsSegundo.Range("G16").Value = 92.09999999999999
$wsSegundo.Range("H16").Value = 7.9
$wsSegundo.Range("I16").Value = 8.5
$wsSegundo.Range("J16").Value = 0
$wsSegundo.Range("K16").Value = 0
$wsSegundo.Range("E17").Value = 346
$wsSegundo.Range("F17").Value = 30
$wsSegundo.Range("G17").Value = 92
$wsSegundo.Range("H17").Value = 8
$wsSegundo.Range("I17").Value = 7.9
$wsSegundo.Range("J17").Value = 0
$wsSegundo.Range("K17").Value = 0
$wsSegundo.Range("E18").Value = 21
$wsSegundo.Range("F18").Value = 15
$wsSegundo.Range("G18").Value = 58.3
$wsSegundo.Range("H18").Value = 41.7
$wsSegundo.Range("I18").Value = 6.5
$wsSegundo.Range("J18").Value = 0
$wsSegundo.Range("K18").Value = 0
$wsSegundo.Range("E19").Value = 20
$wsSegundo.Range("F19").Value = 8
$wsSegundo.Range("G19").Value = 71.40000000000001
$wsSegundo.Range("H19").Value = 28.6
$wsSegundo.Range("I19").Value = 6.4
$wsSegundo.Range("J19").Value = 0
$wsSegundo.Range("K19").Value = 0
$wsSegundo.Range("E20").Value = 14
$wsSegundo.Range("F20").Value = 9
$wsSegundo.Range("G20").Value = 60.9
$wsSegundo.Range("H20").Value = 39.1
$wsSegundo.Range("I20").Value = 6.4
$wsSegundo.Range("J20").Value = 1
$wsSegundo.Range("K20").Value = 4.35
$wsSegundo.Range("E21").Value = 20
$wsSegundo.Range("F21").Value = 9
$wsSegundo.Range("G21").Value = 69
$wsSegundo.Range("H21").Value = 31
$wsSegundo.Range("I21").Value = 6.6
$wsSegundo.Range("J21").Value = 0
$wsSegundo.Range("K21").Value = 0
$wsSegundo.Range("E22").Value = 9
$wsSegundo.Range("F22").Value = 2
$wsSegundo.Range("G22").Value = 81.8
$wsSegundo.Range("H22").Value = 18.2
$wsSegundo.Range("I22").Value = 6.8
$wsSegundo.Range("J22").Value = 0
$wsSegundo.Range("K22").Value = 0
$wsSegundo.Range("E23").Value = 16
$wsSegundo.Range("F23").Value = 3
$wsSegundo.Range("G23").Value = 84.2
$wsSegundo.Range("H23").Value = 15.8
$wsSegundo.Range("I23").Value = 6.1
$wsSegundo.Range("J23").Value = 0
$wsSegundo.Range("K23").Value = 0
$wsSegundo.Range("E24").Value = 16
$wsSegundo.Range("F24").Value = 9
$wsSegundo.Range("G24").Value = 64
$wsSegundo.Range("H24").Value = 36
$wsSegundo.Range("I24").Value = 6.5
$wsSegundo.Range("J24").Value = 0
$wsSegundo.Range("K24").Value = 0
$wsSegundo.Range("E25").Value = 8
$wsSegundo.Range("F25").Value = 5
$wsSegundo.Range("G25").Value = 61.5
$wsSegundo.Range("H25").Value = 38.5
$wsSegundo.Range("I25").Value = 6.2
$wsSegundo.Range("J25").Value = 0
$wsSegundo.Range("K25").Value = 0
$wsSegundo.Range("E26").Value = 13
$wsSegundo.Range("F26").Value = 1
$wsSegundo.Range("G26").Value = 92.90000000000001
$wsSegundo.Range("H26").Value = 7.1
$wsSegundo.Range("I26").Value = 6.7
$wsSegundo.Range("J26").Value = 0
$wsSegundo.Range("K26").Value = 0
$wsSegundo.Range("E27").Value = 137
$wsSegundo.Range("F27").Value = 61
$wsSegundo.Range("G27").Value = 69.2
$wsSegundo.Range("H27").Value = 30.8
$wsSegundo.Range("I27").Value = 6.5
$wsSegundo.Range("J27").Value = 1
$wsSegundo.Range("K27").Value = 0.5
$wsSegundo.Range("E28").Value = 493
$wsSegundo.Range("F28").Value = 140
$wsSegundo.Range("G28").Value = 77.90000000000001
$wsSegundo.Range("H28").Value = 22.1
$wsSegundo.Range("I28").Value = 7
$wsSegundo.Range("J28").Value = 50
$wsSegundo.Range("K28").Value = 0

# --- "Final" sheet updates ---
$wsFinal = $wb.Worksheets.Item("Final")
$wsFinal.Range("I2").Value = 9.5
$wsFinal.Range("I3").Value = 9.5
$wsFinal.Range("E6").Value = 29
$wsFinal.Range("F6").Value = 8
$wsFinal.Range("G6").Value = 78.40000000000001
$wsFinal.Range("H6").Value = 21.6
$wsFinal.Range("I7").Value = 8.699999999999999
$wsFinal.Range("E8").Value = 27
$wsFinal.Range("F8").Value = 5
$wsFinal.Range("G8").Value = 84.40000000000001
$wsFinal.Range("H8").Value = 15.6
$wsFinal.Range("I8").Value = 7.9
$wsFinal.Range("E9").Value = 35
$wsFinal.Range("F9").Value = 4
$wsFinal.Range("G9").Value = 89.7
$wsFinal.Range("H9").Value = 10.3
$wsFinal.Range("I9").Value = 7.8
$wsFinal.Range("E10").Value = 32
$wsFinal.Range("F10").Value = 6
$wsFinal.Range("G10").Value = 84.2
$wsFinal.Range("H10").Value = 15.8
$wsFinal.Range("I10").Value = 8
$wsFinal.Range("I11").Value = 8.1
$wsFinal.Range("I12").Value = 8.6
$wsFinal.Range("I13").Value = 8.800000000000001
$wsFinal.Range("I14").Value = 9.300000000000001
$wsFinal.Range("E16").Value = 35
$wsFinal.Range("F16").Value = 3
$wsFinal.Range("G16").Value = 92.09999999999999
$wsFinal.Range("H16").Value = 7.9
$wsFinal.Range("I16").Value = 9
$wsFinal.Range("E17").Value = 346
$wsFinal.Range("F17").Value = 30
$wsFinal.Range("G17").Value = 92
$wsFinal.Range("H17").Value = 8
$wsFinal.Range("I17").Value = 8.4
$wsFinal.Range("E18").Value = 21
$wsFinal.Range("F18").Value = 15
$wsFinal.Range("G18").Value = 58.3
$wsFinal.Range("H18").Value = 41.7
$wsFinal.Range("I18").Value = 6
$wsFinal.Range("I19").Value = 6.1
$wsFinal.Range("E20").Value = 14
$wsFinal.Range("F20").Value = 9
$wsFinal.Range("G20").Value = 60.9
$wsFinal.Range("H20").Value = 39.1
$wsFinal.Range("I20").Value = 6.1
$wsFinal.Range("E21").Value = 20
$wsFinal.Range("F21").Value = 9
$wsFinal.Range("G21").Value = 69
$wsFinal.Range("H21").Value = 31
$wsFinal.Range("I21").Value = 6.2
$wsFinal.Range("E22").Value = 9
$wsFinal.Range("F22").Value = 2
$wsFinal.Range("G22").Value = 81.8
$wsFinal.Range("H22").Value = 18.2
$wsFinal.Range("I22").Value = 6.6
$wsFinal.Range("E24").Value = 16
$wsFinal.Range("F24").Value = 9
$wsFinal.Range("G24").Value = 64
$wsFinal.Range("H24").Value = 36
$wsFinal.Range("I24").Value = 6
$wsFinal.Range("I25").Value = 6.3
$wsFinal.Range("I26").Value = 6.6
$wsFinal.Range("E27").Value = 137
$wsFinal.Range("F27").Value = 61
$wsFinal.Range("G27").Value = 69.2
$wsFinal.Range("H27").Value = 30.8
$wsFinal.Range("I27").Value = 6.2
$wsFinal.Range("E28").Value = 541
$wsFinal.Range("F28").Value = 92
$wsFinal.Range("G28").Value = 85.5
$wsFinal.Range("H28").Value = 14.5
